$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 33310342
$ws.Range("I33").Value = 49965290
$ws.Range("K33").Value = 49965290
$ws.Range("M33").Value = -49965061

$ws.Range("H129").Value = 982.971
$ws.Range("J129").Value = 1046.129
$ws.Range("L129").Value = 3138.387
$ws.Range("N129").Value = -13138.387

$ws.Range("H132").Value = 833.24243
$ws.Range("I132").Value = 757.7742
$ws.Range("J132").Value = 2003
$ws.Range("K132").Value = 2273.3226
$ws.Range("L132").Value = 6009
$ws.Range("M132").Value = 256.6774
$ws.Range("N132").Value = -11069

$ws.Range("H138").Value = 3611.2954
$ws.Range("I138").Value = 3395.2222
$ws.Range("J138").Value = 3666.8572
$ws.Range("K138").Value = 10185.6666
$ws.Range("L138").Value = 11000.5716
$ws.Range("M138").Value = -5045.6666
$ws.Range("N138").Value = -21280.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2356
$ws.Range("I2").Value = 1545.5
$ws.Range("J2").Value = 3166.5
$ws.Range("K2").Value = 1545.5
$ws.Range("L2").Value = 3166.5
$ws.Range("M2").Value = -1432.5
$ws.Range("N2").Value = -3392.5

$ws.Range("H74").Value = 5380.0557
$ws.Range("I74").Value = 2516.8572
$ws.Range("J74").Value = 9388.532999999999
$ws.Range("K74").Value = 2516.8572
$ws.Range("L74").Value = 9388.532999999999
$ws.Range("M74").Value = -1642.8572
$ws.Range("N74").Value = -11136.533

$ws.Range("H77").Value = 5380.0557
$ws.Range("I77").Value = 2516.8572
$ws.Range("J77").Value = 9388.532999999999
$ws.Range("K77").Value = 12584.286
$ws.Range("L77").Value = 46942.66499999999
$ws.Range("M77").Value = -8216.286
$ws.Range("N77").Value = -55678.66499999999

$ws.Range("H116").Value = 2356
$ws.Range("I116").Value = 1545.5
$ws.Range("J116").Value = 3166.5
$ws.Range("K116").Value = 1545.5
$ws.Range("L116").Value = 3166.5
$ws.Range("M116").Value = 748.5
$ws.Range("N116").Value = -7754.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2356
$ws.Range("I3").Value = 1545.5
$ws.Range("J3").Value = 3166.5
$ws.Range("K3").Value = 1545.5
$ws.Range("L3").Value = 3166.5
$ws.Range("M3").Value = -1431.5
$ws.Range("N3").Value = -3394.5

$ws.Range("H86").Value = 1960.9855
$ws.Range("I86").Value = 1717.2069
$ws.Range("J86").Value = 3246.3635
$ws.Range("K86").Value = 1717.2069
$ws.Range("L86").Value = 3246.3635
$ws.Range("M86").Value = -594.2068999999999
$ws.Range("N86").Value = -5492.363499999999

$ws.Range("H89").Value = 1960.9855
$ws.Range("I89").Value = 1717.2069
$ws.Range("J89").Value = 3246.3635
$ws.Range("K89").Value = 8586.0345
$ws.Range("L89").Value = 16231.8175
$ws.Range("M89").Value = -2970.0345
$ws.Range("N89").Value = -27463.8175

$ws.Range("H138").Value = 67998.336
$ws.Range("J138").Value = 67998.336
$ws.Range("L138").Value = 67998.336
$ws.Range("N138").Value = -78278.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 610185.1
$ws.Range("I31").Value = 4534.963
$ws.Range("J31").Value = 1091142.6
$ws.Range("K31").Value = 4534.963
$ws.Range("L31").Value = 1091142.6
$ws.Range("M31").Value = -4239.963
$ws.Range("N31").Value = -1091732.6

$ws.Range("H34").Value = 610185.1
$ws.Range("I34").Value = 4534.963
$ws.Range("J34").Value = 1091142.6
$ws.Range("K34").Value = 4534.963
$ws.Range("L34").Value = 1091142.6
$ws.Range("M34").Value = -4332.963
$ws.Range("N34").Value = -1091546.6

$ws.Range("H58").Value = 2118840.5
$ws.Range("I58").Value = 3248192
$ws.Range("J58").Value = 10717.6
$ws.Range("K58").Value = 3248192
$ws.Range("L58").Value = 10717.6
$ws.Range("M58").Value = -3247989
$ws.Range("N58").Value = -11123.6

$ws.Range("H99").Value = 1125
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 1166.6666
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 1166.6666
$ws.Range("M99").Value = 498
$ws.Range("N99").Value = -4162.6666

$ws.Range("H126").Value = 1125
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 1166.6666
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 3499.9998
$ws.Range("M126").Value = -530
$ws.Range("N126").Value = -8439.9998

$ws.Range("H132").Value = 4193.6665
$ws.Range("I132").Value = 4355.2354
$ws.Range("J132").Value = 3507
$ws.Range("K132").Value = 13065.7062
$ws.Range("L132").Value = 10521
$ws.Range("M132").Value = -10535.7062
$ws.Range("N132").Value = -15581

$ws.Range("H134").Value = 3383.2273
$ws.Range("I134").Value = 2962.6875
$ws.Range("J134").Value = 4504.6665
$ws.Range("K134").Value = 8888.0625
$ws.Range("L134").Value = 13513.9995
$ws.Range("M134").Value = -6353.0625
$ws.Range("N134").Value = -18583.9995

$ws.Range("H136").Value = 2118840.5
$ws.Range("I136").Value = 3248192
$ws.Range("J136").Value = 10717.6
$ws.Range("K136").Value = 9744576
$ws.Range("L136").Value = 32152.8
$ws.Range("M136").Value = -9742026
$ws.Range("N136").Value = -37252.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2631.0154
$ws.Range("J68").Value = 4374.4243
$ws.Range("L68").Value = 13123.2729
$ws.Range("N68").Value = -14745.2729

$ws.Range("H71").Value = 2631.0154
$ws.Range("J71").Value = 4374.4243
$ws.Range("L71").Value = 39369.8187
$ws.Range("N71").Value = -47481.8187

$ws.Range("H113").Value = 550.5143
$ws.Range("I113").Value = 536.0702
$ws.Range("J113").Value = 613.8461
$ws.Range("K113").Value = 1608.2106
$ws.Range("L113").Value = 1841.5383
$ws.Range("M113").Value = 561.7894000000001
$ws.Range("N113").Value = -6181.5383

$ws.Range("H131").Value = 36547.035
$ws.Range("I131").Value = 1211.4286
$ws.Range("J131").Value = 48914.5
$ws.Range("K131").Value = 3634.2858
$ws.Range("L131").Value = 146743.5
$ws.Range("M131").Value = 1405.7142
$ws.Range("N131").Value = -156823.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3548.44
$ws.Range("I102").Value = 2957.8235
$ws.Range("J102").Value = 4803.5
$ws.Range("K102").Value = 2957.8235
$ws.Range("L102").Value = 4803.5
$ws.Range("M102").Value = -1335.8235
$ws.Range("N102").Value = -8047.5

$ws.Range("H122").Value = 11792.615
$ws.Range("I122").Value = 13859.6
$ws.Range("J122").Value = 4902.6665
$ws.Range("K122").Value = 41578.8
$ws.Range("L122").Value = 14707.9995
$ws.Range("M122").Value = -39128.8
$ws.Range("N122").Value = -19607.9995

$ws.Range("H134").Value = 44268.4
$ws.Range("J134").Value = 44268.4
$ws.Range("L134").Value = 132805.2
$ws.Range("N134").Value = -137875.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2800.8
$ws.Range("I7").Value = 2251
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 2251
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -2139
$ws.Range("N7").Value = -5224

$ws.Range("H22").Value = 4062.625
$ws.Range("J22").Value = 4571.5713
$ws.Range("L22").Value = 4571.5713
$ws.Range("N22").Value = -5161.5713

$ws.Range("H27").Value = 4062.625
$ws.Range("J27").Value = 4571.5713
$ws.Range("L27").Value = 4571.5713
$ws.Range("N27").Value = -4785.5713

$ws.Range("H40").Value = 4386.2856
$ws.Range("I40").Value = 3450.6667
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 3450.6667
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -3314.6667
$ws.Range("N40").Value = -10272

$ws.Range("H46").Value = 1033.3334
$ws.Range("I46").Value = 771.4286
$ws.Range("K46").Value = 771.4286
$ws.Range("M46").Value = -583.4286

$ws.Range("H61").Value = 40125.5
$ws.Range("I61").Value = 40125.5
$ws.Range("K61").Value = 40125.5
$ws.Range("M61").Value = -39923.5

$ws.Range("H113").Value = 40125.5
$ws.Range("I113").Value = 40125.5
$ws.Range("K113").Value = 40125.5
$ws.Range("M113").Value = -37955.5

$ws.Range("H122").Value = 5665.2812
$ws.Range("I122").Value = 5617.2856
$ws.Range("J122").Value = 6001.25
$ws.Range("K122").Value = 16851.8568
$ws.Range("L122").Value = 18003.75
$ws.Range("M122").Value = -14401.8568
$ws.Range("N122").Value = -22903.75

$ws.Range("H126").Value = 2800.8
$ws.Range("I126").Value = 2251
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 6753
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -4283
$ws.Range("N126").Value = -19940

$ws.Range("H135").Value = 75448
$ws.Range("J135").Value = 75448
$ws.Range("L135").Value = 75448
$ws.Range("N135").Value = -85588

$ws.Range("H141").Value = 79861.664
$ws.Range("J141").Value = 79861.664
$ws.Range("L141").Value = 79861.664
$ws.Range("N141").Value = -90221.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 872.8421
$ws.Range("I113").Value = 493.33334
$ws.Range("J113").Value = 944
$ws.Range("K113").Value = 1480.00002
$ws.Range("L113").Value = 2832
$ws.Range("M113").Value = 689.9999800000001
$ws.Range("N113").Value = -7172

$ws.Range("H133").Value = 53614
$ws.Range("J133").Value = 53614
$ws.Range("L133").Value = 53614
$ws.Range("N133").Value = -63734

$ws.Range("H136").Value = 6202.7295
$ws.Range("I136").Value = 2361.8572
$ws.Range("J136").Value = 8540.652
$ws.Range("K136").Value = 7085.571599999999
$ws.Range("L136").Value = 25621.956
$ws.Range("M136").Value = -4535.571599999999
$ws.Range("N136").Value = -30721.956

$ws.Range("H141").Value = 39833.75
$ws.Range("J141").Value = 39833.75
$ws.Range("L141").Value = 39833.75
$ws.Range("N141").Value = -50193.75
